$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from H1, the last header cell) onto the
# two new header cells so they match the bold/bordered header formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-58
$values = @(
    @(2, 6, 6),
    @(3, 4, 5),
    @(4, 6, 6),
    @(5, 6, 7),
    @(6, 7, 7),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 7, 7),
    @(10, 6, 6),
    @(11, 6, 6),
    @(12, 4, 4),
    @(13, 6, 7),
    @(14, 11, 11),
    @(15, 5, 6),
    @(16, 9, 10),
    @(17, 9, 9),
    @(18, 7, 7),
    @(19, 9, 9),
    @(20, 7, 8),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 7, 7),
    @(26, 6, 7),
    @(27, 7, 7),
    @(28, 7, 8),
    @(29, 5, 5),
    @(30, 7, 7),
    @(31, 5, 6),
    @(32, 8, 8),
    @(33, 7, 7),
    @(34, 8, 8),
    @(35, 9, 9),
    @(36, 10, 10),
    @(37, 6, 7),
    @(38, 9, 9),
    @(39, 6, 6),
    @(40, 8, 8),
    @(41, 6, 7),
    @(42, 7, 8),
    @(43, 8, 8),
    @(44, 5, 6),
    @(45, 12, 13),
    @(46, 7, 7),
    @(47, 6, 7),
    @(48, 6, 6),
    @(49, 6, 6),
    @(50, 7, 7),
    @(51, 7, 8),
    @(52, 8, 8),
    @(53, 7, 7),
    @(54, 6, 7),
    @(55, 6, 6),
    @(56, 7, 8),
    @(57, 8, 8),
    @(58, 4, 4)
)

foreach ($row in $values) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}
